$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch area used as temporary holding space while swapping rows
$scratch = $ws.Range("Z1:AD1")

# For every year block (rows base..base+3 = quarters A,B,C,D), the "B" (quarter 2)
# and "C" (quarter 3) rows need their A:E content swapped.
for ($base = 2; $base -le 62; $base += 4) {
    $rB = $base + 1
    $rC = $base + 2

    $rangeB = $ws.Range("A" + $rB + ":E" + $rB)
    $rangeC = $ws.Range("A" + $rC + ":E" + $rC)

    $rangeB.Copy($scratch)
    $rangeC.Copy($rangeB)
    $scratch.Copy($rangeC)
}
$scratch.Clear()

# Columns F (产销率) and G (销售量) were removed entirely.
$ws.Range("F1:G65").Delete()
